# Fruta / hortaliza, semanal
# Insert the latest weekly price-report row for
# "Feria Lagunitas de Puerto Montt - Espinaca" at the top of the data
# block (row 16, right after the most-recent existing entry), pushing
# every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16:71 down to 17:72, creating a blank row 16.
$ws.Rows.Item(16).Insert()

# Populate the new row with this week's figures.
$ws.Range("A16").Value2 = 4
$ws.Range("B16").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C16").Value2 = "Los Lagos"
$ws.Range("D16").Value2 = 45107
$ws.Range("E16").Value2 = 10
$ws.Range("F16").Value2 = 100112012
$ws.Range("G16").Value2 = "Espinaca"
$ws.Range("H16").Value2 = "Sin especificar"
$ws.Range("I16").Value2 = "Primera"
$ws.Range("J16").Value2 = 25
$ws.Range("K16").Value2 = 14000
$ws.Range("L16").Value2 = 14000
$ws.Range("M16").Value2 = 14000
$ws.Range("N16").Value2 = "`$/cuna 10 kilos"
$ws.Range("O16").Value2 = "Región Metropolitana"
$ws.Range("P16").Value2 = 1400
$ws.Range("Q16").Value2 = 10
$ws.Range("R16").Value2 = "Hortaliza"
